$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new column D
$ws.Columns("D").Insert()

# Unlock D1 so we can write to it (sheet is protected)
$ws.Range("D1").Locked = $false
$ws.Range("D1").Value = "org_email"

# Restore D1's original header style (center+wrap+bold-ish font, locked) by
# copying formatting from the C1 header cell.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# Re-protect (drops the password hash, same effect as re-saving protection)
$ws.Protect()

# Move selection to match target
$ws.Range("F7").Select()
Write-Host "Done"
